# Update "through" date from April 18 to April 19, 2022, and add the
# carjacking counts that came in for that extra day across the various
# "April" (month-to-date) columns for each year.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename sheet & title text ------------------------------------------------
$ws.Name = "Through 2022-04-19"
$ws.Range("B1").Value = "April 2022 (through April 19)"

# --- Row 2: Austin -------------------------------------------------------------
$ws.Range("J2").Value = 8   # was 7
$ws.Range("N2").Value = 5   # was 4
$ws.Range("Z2").Value = 1   # new

# --- Row 3: Englewood -----------------------------------------------------------
$ws.Range("Z3").Value = 2   # was 1

# --- Row 4: North Lawndale -------------------------------------------------------
$ws.Range("F4").Value = 6   # was 5

# --- Row 15: Washington Park -------------------------------------------------
$ws.Range("F15").Value = 3  # was 2

# --- Row 24: Lincoln Park ------------------------------------------------------
$ws.Range("B24").Value = 2  # was 1

# --- Row 27: Uptown --------------------------------------------------------------
$ws.Range("R27").Value = 1  # new

# --- Row 28: West Loop -----------------------------------------------------------
$ws.Range("V28").Value = 1  # new

# --- Row 39: Little Village --------------------------------------------------
$ws.Range("B39").Value = 1  # new
$ws.Range("R39").Value = 1  # new

# --- Row 44: Brighton Park ---------------------------------------------------
$ws.Range("N44").Value = 1  # new

# --- Row 50: Grand Crossing --------------------------------------------------
$ws.Range("J50").Value = 3  # was 2

# --- Row 60: Ashburn -------------------------------------------------------------
$ws.Range("V60").Value = 2  # was 1

# --- Row 68: Fuller Park ------------------------------------------------------
$ws.Range("B68").Value = 1  # new

# --- Row 75: Little Italy, UIC -----------------------------------------------
$ws.Range("B75").Value = 3  # was 2
